$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the value into A1 (matches <row r="1"><c r="A1"><v>123123</v></c></row>)
$ws.Range("A1").Value = 123123

# Move / leave the active selection on B8 (matches <selection activeCell="B8" sqref="B8"/>)
$ws.Range("B8").Select() | Out-Null
